$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.898.07'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.874.46'
$ws.Range('E3').Value = '  -1.02%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9982'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.24%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.7387'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -4.94%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '242.02'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.73%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.9997'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3158'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.83%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07193'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.65%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '24.73'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -4.23%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.08359'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -4.08%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.7505'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.414'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '1.857.32'
$ws.Range('E14').Value = '  -11.36%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '92.59'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').Value = '29.898.77'
$ws.Range('E16').Value = '  -1.47%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '6.070'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.93%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '246.33'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000007835'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.27%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.9986'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = '2.125.70'
$ws.Range('E22').Value = '  -11.77%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '8.000'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.9992'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -6.11%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.261'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.32%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '165.19'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.17%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '18.67'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.90%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.033'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  +5.13%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.580'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.45%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.536'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.47%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.272'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +3.17%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.05315'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.89%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.237'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.62%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.7544'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.14%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.9991'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.06%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.698'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.26%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01957'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('E40').Value = '  -1.18%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.4512'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = '1.113.20'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.050'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.20%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '72.29'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.79%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.8545'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '104.40'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.04%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.630'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '3.090'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.61%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.841'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').Value = '2.022.87'
$ws.Range('E51').Value = '  -10.18%  '
